$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "A 14042-2023") now holds the data previously in row 3 ("A 5398-2026"),
# and row 3 now holds what used to be row 2's data. Rows 5/6/7 are cycled and rows
# 8/9 are swapped in the same manner. All "Förändrad" (column C) dates are bumped
# to 46070 for every data row.
$ws.Range("A2").Value = 'A 5398-2026'
$ws.Range("B2").Value = 46050.49048611111
$ws.Range("C2").Value = 46070
$ws.Range("G2").Value = 0.5
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2
$ws.Range("R2").Value = "Brandticka`r`nKambräken"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/artfynd/A 5398-2026 artfynd.xlsx", "A 5398-2026")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/kartor/A 5398-2026 karta.png", "A 5398-2026")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomål/A 5398-2026 FSC-klagomål.docx", "A 5398-2026")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomålsmail/A 5398-2026 FSC-klagomål mail.docx", "A 5398-2026")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsyn/A 5398-2026 tillsynsbegäran.docx", "A 5398-2026")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsynsmail/A 5398-2026 tillsynsbegäran mail.docx", "A 5398-2026")'
$ws.Range("A3").Value = 'A 14042-2023'
$ws.Range("B3").Value = 45008
$ws.Range("C3").Value = 46070
$ws.Range("G3").Value = 4.1
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 0
$ws.Range("R3").Value = "Revlummer`r`nÄkta lopplummer"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/artfynd/A 14042-2023 artfynd.xlsx", "A 14042-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/kartor/A 14042-2023 karta.png", "A 14042-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomål/A 14042-2023 FSC-klagomål.docx", "A 14042-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomålsmail/A 14042-2023 FSC-klagomål mail.docx", "A 14042-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsyn/A 14042-2023 tillsynsbegäran.docx", "A 14042-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsynsmail/A 14042-2023 tillsynsbegäran mail.docx", "A 14042-2023")'
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Range("C4").Value = 46070
$ws.Range("A5").Value = 'A 25251-2025'
$ws.Range("B5").Value = 45800.50082175926
$ws.Range("C5").Value = 46070
$ws.Range("G5").Value = 0.7
$ws.Range("A6").Value = 'A 25254-2025'
$ws.Range("B6").Value = 45800.50479166667
$ws.Range("C6").Value = 46070
$ws.Range("F6").Value = 'Kommuner'
$ws.Range("G6").Value = 0.2
$ws.Range("A7").Value = 'A 26074-2025'
$ws.Range("B7").Value = 45805.32366898148
$ws.Range("C7").Value = 46070
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 1.3
$ws.Range("A8").Value = 'A 5402-2026'
$ws.Range("B8").Value = 46050.49721064815
$ws.Range("C8").Value = 46070
$ws.Range("F8").Value = 'Kommuner'
$ws.Range("G8").Value = 0.7
$ws.Range("A9").Value = 'A 35734-2023'
$ws.Range("B9").Value = 45147.89258101852
$ws.Range("C9").Value = 46070
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 5.9
